# CI: Update Excel counters (state_counters + packages)
# Appends newly-tracked Maryland SPA/Waiver packages to the "Packages" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=PackageType, B=State, C=Authority, D=ActionType, E=PackageID, F=Status, G=ParentID
$newRows = @(
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9423",     "Under Review",       ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9424",     "Approved",           ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9425",     "Submitted",          ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9426",     "Under Review",       ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.34", "Pending-Approval",   "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9427",     "Pending-Concurrence",""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9428",     "Submitted",          ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9429",     "Submitted",          ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.35", "Unsubmitted",        "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9430",     "Under Review",       ""),
    @("Waiver", "MD", "1915(b)",      "Initial",   "MD-2275.R00.00", "Terminated",         ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9431",     "Under Review",       "")
)

$startRow = $ws.UsedRange.Rows.Count + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
